# Actualización automática 2025-12-11 14:30:07
$wb = $excel.ActiveWorkbook

# --- Hoja "VENTAS POR GRUPO" ---
$wsGrupo = $wb.Worksheets.Item("VENTAS POR GRUPO")

# Fila 11 (ILLER LOPEZ ROBERTO FERNANDO / FRANK FERRETERIA FRANKFERRE CIA.)
$wsGrupo.Range("E11").Value = 84.59
$wsGrupo.Range("I11").Value = 180
$wsGrupo.Range("M11").Value = 4610.35

# Fila 21 (totales "X de 19")
$wsGrupo.Range("E21").Value = "3 de 19"
$wsGrupo.Range("I21").Value = "1 de 19"
$wsGrupo.Range("M21").Value = "2 de 19"

# --- Hoja "VENTA MENSUAL" ---
$wsMensual = $wb.Worksheets.Item("VENTA MENSUAL")

# Fila 11 (ILLER LOPEZ ROBERTO FERNANDO / FRANK FERRETERIA FRANKFERRE CIA.)
$wsMensual.Range("F11").Value = 4874.94

# Fila 21 (totales)
$wsMensual.Range("F21").Value = 7340.91
